$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Maria Claudia%Alzamora%NULL%0,                    Tania%Paredes%NULL%0,                    Tania%Paredes%NULL%0,                    David%Caceres%NULL%0,                    Camille M.%Webb%NULL%0,                    Luis M.%Valdez%NULL%0,                    Mauricio%La Rosa%NULL%0]'
$ws.Range("I2").Value = '_PMC_CROSSREF'

$ws.Range("E3").Value = '[Rachel M.%Burke%NULL%1,                    Claire M.%Midgley%NULL%1,                    Alissa%Dratch%NULL%1,                    Marty%Fenstersheib%NULL%1,                    Thomas%Haupt%NULL%1,                    Michelle%Holshue%NULL%1,                    Isaac%Ghinai%NULL%1,                    M. Claire%Jarashow%NULL%1,                    Jennifer%Lo%NULL%1,                    Tristan D.%McPherson%NULL%1,                    Sara%Rudman%NULL%1,                    Sarah%Scott%NULL%1,                    Aron J.%Hall%NULL%1,                    Alicia M.%Fry%NULL%1,                    Melissa A.%Rolfes%NULL%1]'
$ws.Range("I3").Value = '_PMC_CROSSREF'

$ws.Range("E4").Value = '[Jing%Cai%NULL%1,                    Wenjie%Sun%NULL%1,                    Jianping%Huang%NULL%1,                    Michelle%Gamber%NULL%1,                    Jing%Wu%NULL%2,                    Guiqing%He%NULL%1]'
$ws.Range("I4").Value = '_PMC_CROSSREF'

$ws.Range("E5").Value = '[Jasper Fuk-Woo%Chan%NULL%0,                    Shuofeng%Yuan%NULL%1,                    Kin-Hang%Kok%NULL%1,                    Kelvin Kai-Wang%To%NULL%1,                    Hin%Chu%NULL%1,                    Jin%Yang%NULL%1,                    Fanfan%Xing%NULL%1,                    Jieling%Liu%NULL%1,                    Cyril Chik-Yan%Yip%NULL%2,                    Rosana Wing-Shan%Poon%NULL%2,                    Hoi-Wah%Tsoi%NULL%1,                    Simon Kam-Fai%Lo%NULL%1,                    Kwok-Hung%Chan%NULL%2,                    Vincent Kwok-Man%Poon%NULL%1,                    Wan-Mui%Chan%NULL%2,                    Jonathan Daniel%Ip%NULL%2,                    Jian-Piao%Cai%NULL%2,                    Vincent Chi-Chung%Cheng%NULL%2,                    Honglin%Chen%NULL%2,                    Christopher Kim-Ming%Hui%NULL%1,                    Kwok-Yung%Yuen%NULL%0]'
$ws.Range("I5").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E6").Value = '[Huijun%Chen%NULL%0,                    Juanjuan%Guo%NULL%0,                    Chen%Wang%NULL%0,                    Fan%Luo%NULL%1,                    Xuechen%Yu%NULL%0,                    Wei%Zhang%NULL%0,                    Jiafu%Li%NULL%0,                    Dongchi%Zhao%NULL%1,                    Dan%Xu%NULL%3,                    Qing%Gong%NULL%1,                    Jing%Liao%NULL%1,                    Huixia%Yang%yanghuixia@bjmu.edu.cn%0,                    Wei%Hou%houwei@whu.edu.cn%0,                    Yuanzhen%Zhang%zhangyuanzhen@whu.edu.cn%0]'
$ws.Range("I6").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E7").Value = '[Yifei%Chen%xref no email%1,            Liangjun%Chen%xref no email%1,            Qiaoling%Deng%xref no email%0,            Guqin%Zhang%xref no email%1,            Kaisong%Wu%xref no email%1,            Lan%Ni%xref no email%1,            Yibin%Yang%xref no email%1,            Bing%Liu%xref no email%1,            Wei%Wang%xref no email%0,            Chaojie%Wei%xref no email%1,            Jiong%Yang%xref no email%1,            Guangming%Ye%xref no email%1,            Zhenshun%Cheng%xref no email%1]'
$ws.Range("I7").Value = '_CROSSREF'

$ws.Range("E8").Value = '[Yan%Chen%NULL%2,                    Hua%Peng%NULL%1,                    Lin%Wang%NULL%3,                    Yin%Zhao%NULL%1,                    Lingkong%Zeng%NULL%1,                    Hui%Gao%NULL%1,                    Yalan%Liu%NULL%1]'
$ws.Range("I8").Value = '_PMC_CROSSREF'

$ws.Range("E9").Value = '[Hao-Yuan%Cheng%NULL%1,                    Shu-Wan%Jian%NULL%1,                    Ding-Ping%Liu%NULL%1,                    Ta-Chou%Ng%NULL%1,                    Wan-Ting%Huang%NULL%1,                    Hsien-Ho%Lin%NULL%1]'
$ws.Range("I9").Value = '_PMC_CROSSREF'

$ws.Range("E10").Value = '[Vincent C. C.%Cheng%NULL%1,                    Shuk-Ching%Wong%NULL%1,                    Jonathan H. K.%Chen%NULL%1,                    Cyril C. Y.%Yip%NULL%1,                    Vivien W. M.%Chuang%NULL%1,                    Owen T. Y.%Tsang%NULL%1,                    Siddharth%Sridhar%NULL%1,                    Jasper F. W.%Chan%NULL%1,                    Pak-Leung%Ho%NULL%1,                    Kwok-Yung%Yuen%NULL%0]'
$ws.Range("I10").Value = '_PMC_CROSSREF'

$ws.Range("E11").Value = '[Po Ying%Chia%NULL%1,                    Kristen Kelli%Coleman%NULL%2,                    Kristen Kelli%Coleman%NULL%0,                    Yian Kim%Tan%NULL%2,                    Yian Kim%Tan%NULL%0,                    Sean Wei Xiang%Ong%NULL%1,                    Marcus%Gum%NULL%2,                    Marcus%Gum%NULL%0,                    Sok Kiang%Lau%NULL%1,                    Xiao Fang%Lim%NULL%1,                    Ai Sim%Lim%NULL%1,                    Stephanie%Sutjipto%NULL%1,                    Pei Hua%Lee%NULL%1,                    Than The%Son%NULL%1,                    Barnaby Edward%Young%NULL%1,                    Donald K.%Milton%NULL%1,                    Gregory C.%Gray%NULL%2,                    Gregory C.%Gray%NULL%0,                    Stephan%Schuster%NULL%2,                    Stephan%Schuster%NULL%0,                    Timothy%Barkham%NULL%1,                    Partha Pratim%De%NULL%2,                    Partha Pratim%De%NULL%0,                    Shawn%Vasoo%NULL%1,                    Monica%Chan%NULL%1,                    Brenda Sze Peng%Ang%NULL%1,                    Boon Huan%Tan%NULL%1,                    Yee-Sin%Leo%NULL%1,                    Oon-Tek%Ng%Oon_Tek_Ng@ncid.sg%1,                    Michelle Su Yen%Wong%NULL%2,                    Michelle Su Yen%Wong%NULL%0,                    Kalisvar%Marimuthu%kalisvar_marimuthu@ncid.sg%1,                    David Chien%Lye%NULL%2,                    David Chien%Lye%NULL%0,                    Poh Lian%Lim%NULL%1,                    Cheng Chuan%Lee%NULL%1,                    Li Min%Ling%NULL%1,                    Lawrence%Lee%NULL%1,                    Tau Hong%Lee%NULL%1,                    Chen Seong%Wong%NULL%1,                    Sapna%Sadarangani%NULL%1,                    Ray Junhao%Lin%NULL%1,                    Deborah Hee Ling%Ng%NULL%1,                    Mucheli%Sadasiv%NULL%1,                    Tsin Wen%Yeo%NULL%1,                    Chiaw Yee%Choy%NULL%1,                    Glorijoy Shi En%Tan%NULL%1,                    Frederico%Dimatatac%NULL%1,                    Isais Florante%Santos%NULL%1,                    Chi Jong%Go%NULL%1,                    Yu Kit%Chan%NULL%1,                    Jun Yang%Tay%NULL%1,                    Jackie Yu-Ling%Tan%NULL%1,                    Nihar%Pandit%NULL%1,                    Benjamin Choon Heng%Ho%NULL%1,                    Shehara%Mendis%NULL%1,                    Yuan Yi Constance%Chen%NULL%1,                    Mohammad Yazid%Abdad%NULL%1,                    Daniela%Moses%NULL%1]'
$ws.Range("I11").Value = '_PMC_Springer_CROSSREF'

$ws.Range("E12").Value = '[Alex W H%Chin%NULL%1,                    Julie T S%Chu%NULL%1,                    Mahen R A%Perera%NULL%1,                    Kenrie P Y%Hui%NULL%1,                    Hui-Ling%Yen%NULL%1,                    Michael C W%Chan%NULL%1,                    Malik%Peiris%NULL%1,                    Leo L M%Poon%llmpoon@hku.hk%1]'
$ws.Range("I12").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E13").Value = '[Francesca%Colavita%NULL%2,                    Daniele%Lapa%NULL%2,                    Fabrizio%Carletti%NULL%2,                    Eleonora%Lalle%NULL%2,                    Licia%Bordi%NULL%2,                    Patrizia%Marsella%NULL%2,                    Emanuele%Nicastri%NULL%2,                    Nazario%Bevilacqua%NULL%2,                    Maria Letizia%Giancola%NULL%2,                    Angela%Corpolongo%NULL%2,                    Giuseppe%Ippolito%NULL%2,                    Maria Rosaria%Capobianchi%NULL%2,                    Concetta%Castilletti%NULL%2]'
$ws.Range("I13").Value = '_PMC_CROSSREF'

$ws.Range("E14").Value = '[Francesca%Colavita%NULL%0,                    Daniele%Lapa%NULL%0,                    Fabrizio%Carletti%NULL%0,                    Eleonora%Lalle%NULL%0,                    Licia%Bordi%NULL%0,                    Patrizia%Marsella%NULL%0,                    Emanuele%Nicastri%NULL%0,                    Nazario%Bevilacqua%NULL%0,                    Maria Letizia%Giancola%NULL%0,                    Angela%Corpolongo%NULL%0,                    Giuseppe%Ippolito%NULL%0,                    Maria Rosaria%Capobianchi%NULL%0,                    Concetta%Castilletti%NULL%0]'
$ws.Range("I14").Value = '_PMC_CROSSREF'

$ws.Range("E15").Value = '[ P.%Cui%null%1,                     Z.% Chen%null%1,                     T.% Wang%null%2,                     J.% Dai%null%1,                     J.% Zhang%null%1,                     T.% Ding%null%1,                     J.% Jiang%null%1,                     J.% Liu%null%1,                     C.% Zhang%null%1,                     W.% Shan%null%1,                     S.% Wang%null%1,                     Y.% Rong%null%1,                     J.% Chang%null%1,                     X.% Miao%null%1,                     X.% Ma%null%1,                     S. % Wang%null%1]'
$ws.Range("I15").Value = '_MedBiorxiv_CROSSREF'

$ws.Range("E16").Value = '[Lan%Dong%xref no email%0,            Jinhua%Tian%xref no email%0,            Songming%He%xref no email%0,            Chuchao%Zhu%xref no email%0,            Jian%Wang%xref no email%0,            Chen%Liu%xref no email%0,            Jing%Yang%xref no email%0]'
$ws.Range("I16").Value = '_CROSSREF'

$ws.Range("E17").Value = '[ Yawen%Dong%null%1,          Tsung-Hao%Liu%null%1,          Thomas%Yau%null%1,          Chiun%Hsu%null%1]'
$ws.Range("I17").Value = '_CROSSREF'

$ws.Range("E18").Value = '[Qihong%Fan%NULL%1,                    Yan%Pan%woshipanyan@126.com%1,                    Qingcui%Wu%NULL%1,                    Shan%Liu%NULL%1,                    Xu%Song%NULL%1,                    Zhongguo%Xie%NULL%1,                    Yang%Liu%NULL%1,                    Liang%Zhao%NULL%1,                    Zhonghong%Wang%NULL%1,                    Yifei%Zhang%NULL%1,                    Zuchuang%Wu%NULL%1,                    Lei%Guan%NULL%1,                    Xiaolong%Lv%NULL%1]'
$ws.Range("I18").Value = '_PMC_CROSSREF'

$ws.Range("E19").Value = '[Sasan%Faridi%NULL%1,                    Sadegh%Niazi%NULL%1,                    Kaveh%Sadeghi%NULL%1,                    Kazem%Naddafi%NULL%1,                    Jila%Yavarian%NULL%1,                    Mansour%Shamsipour%NULL%1,                    Nazanin Zahra Shafiei%Jandaghi%NULL%1,                    Khosro%Sadeghniiat%NULL%1,                    Ramin%Nabizadeh%NULL%1,                    Masud%Yunesian%NULL%1,                    Fatemeh%Momeniha%NULL%1,                    Adel%Mokamel%NULL%1,                    Mohammad Sadegh%Hassanvand%NULL%1,                    Talat%MokhtariAzad%NULL%1]'
$ws.Range("I19").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E21").Value = '[Mi Seon%Han%NULL%1,                    Moon-Woo%Seong%NULL%1,                    Eun Young%Heo%NULL%1,                    Ji Hong%Park%NULL%1,                    Namhee%Kim%NULL%1,                    Sue%Shin%NULL%1,                    Sung Im%Cho%NULL%1,                    Sung Sup%Park%NULL%1,                    Eun Hwa%Choi%eunchoi@snu.ac.kr%1]'
$ws.Range("I21").Value = '_PMC_CROSSREF'

$ws.Range("E22").Value = '[Amy%Heinzerling%NULL%1,                    Matthew J.%Stuckey%NULL%1,                    Tara%Scheuer%NULL%1,                    Kerui%Xu%NULL%1,                    Kiran M.%Perkins%NULL%1,                    Heather%Resseger%NULL%1,                    Shelley%Magill%NULL%1,                    Jennifer R.%Verani%NULL%1,                    Seema%Jain%NULL%1,                    Meileen%Acosta%NULL%1,                    Erin%Epson%NULL%1]'
$ws.Range("I22").Value = '_PMC_CROSSREF'

$ws.Range("E23").Value = '[Lei%Huang%huangleizhenting@126.com%3,                    Xiuwen%Zhang%NULL%2,                    Xinyue%Zhang%NULL%2,                    Zhijian%Wei%NULL%2,                    Lingli%Zhang%NULL%2,                    Jingjing%Xu%NULL%2,                    Peipei%Liang%NULL%2,                    Yuanhong%Xu%xyhong1964@163.com%0,                    Chengyuan%Zhang%2396476405@qq.com%2,                    Aman%Xu%amanxu@163.com%2]'
$ws.Range("I23").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E24").Value = '[Lei%Huang%huangleizhenting@126.com%0,                    Xiuwen%Zhang%NULL%0,                    Xinyue%Zhang%NULL%0,                    Zhijian%Wei%NULL%0,                    Lingli%Zhang%NULL%0,                    Jingjing%Xu%NULL%0,                    Peipei%Liang%NULL%0,                    Yuanhong%Xu%xyhong1964@163.com%0,                    Chengyuan%Zhang%2396476405@qq.com%0,                    Aman%Xu%amanxu@163.com%0]'
$ws.Range("I24").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E25").Value = '[Rui%Huang%NULL%0,                    Juan%Xia%NULL%1,                    Yuxin%Chen%NULL%1,                    Chun%Shan%NULL%1,                    Chao%Wu%NULL%2]'
$ws.Range("I25").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E26").Value = '[Sara N.%Iqbal%NULL%0,                    Rachael%Overcash%NULL%0,                    Neggin%Mokhtari%NULL%0,                    Haleema%Saeed%NULL%0,                    Stacey%Gold%NULL%0,                    Tamika%Auguste%NULL%0,                    Muhammad-Usman%Mirza%NULL%0,                    Maria-Elena%Ruiz%NULL%0,                    Joeffrey J.%Chahine%NULL%0,                    Joeffrey J.%Chahine%NULL%0,                    Masashi%Waga%NULL%0,                    Masashi%Waga%NULL%0,                    Glenn%Wortmann%NULL%0]'
$ws.Range("I26").Value = '_PMC_CROSSREF'

$ws.Range("E27").Value = '[Xuejun%Jiang%NULL%1,                    Mei%Luo%NULL%1,                    Zhen%Zou%NULL%1,                    Xu%Wang%875777340@qq.com%1,                    Chengzhi%Chen%chengzhichen@cqmu.edu.cn%1,                    Jingfu%Qiu%jfqiu@126.com%2,                    Jingfu%Qiu%jfqiu@126.com%0]'
$ws.Range("I27").Value = '_PMC_CROSSREF'

$ws.Range("E28").Value = '[Xiao-Lin%Jiang%NULL%1,                    Xiao-Li%Zhang%NULL%1,                    Xiang-Na%Zhao%NULL%1,                    Cun-Bao%Li%NULL%1,                    Jie%Lei%NULL%2,                    Zeng-Qiang%Kou%NULL%1,                    Wen-Kui%Sun%NULL%1,                    Yang%Hang%NULL%1,                    Feng%Gao%NULL%1,                    Sheng-Xiang%Ji%NULL%1,                    Can-Fang%Lin%NULL%1,                    Bo%Pang%NULL%1,                    Ming-Xiao%Yao%NULL%1,                    Benjamin D%Anderson%NULL%1,                    Guo-Lin%Wang%NULL%1,                    Lin%Yao%NULL%1,                    Li-Jun%Duan%NULL%1,                    Dian-Ming%Kang%dmkang66@163.com%1,                    Mai-Juan%Ma%mjma@163.com%1]'
$ws.Range("I28").Value = '_PMC_CROSSREF'

$ws.Range("I29").Value = '_PMC_CROSSREF'

$ws.Range("E30").Value = '[Hai T%Le%NULL%1,                    Lam V%Nguyen%NULL%1,                    Dien M%Tran%NULL%1,                    Hai T%Do%NULL%1,                    Huong T%Tran%NULL%1,                    Yen T%Le%NULL%1,                    Phuc H%Phan%phuc.h.phan@nhp.org.vn%1]'
$ws.Range("I30").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E31").Value = '[Dong Hwan%Lee%NULL%0,                    Jihyang%Lee%NULL%0,                    Jihyang%Lee%NULL%0,                    Eunju%Kim%NULL%0,                    Eunju%Kim%NULL%0,                    Kyeongyoon%Woo%NULL%0,                    Kyeongyoon%Woo%NULL%0,                    Hak Youle%Park%NULL%0,                    Hak Youle%Park%NULL%0,                    Jihyun%An%NULL%0,                    Jihyun%An%NULL%0]'
$ws.Range("I31").Value = '_PMC_CROSSREF'

$ws.Range("E32").Value = '[Dong Hwan%Lee%NULL%0,                    Jihyang%Lee%NULL%0,                    Jihyang%Lee%NULL%0,                    Eunju%Kim%NULL%0,                    Eunju%Kim%NULL%0,                    Kyeongyoon%Woo%NULL%0,                    Kyeongyoon%Woo%NULL%0,                    Hak Youle%Park%NULL%0,                    Hak Youle%Park%NULL%0,                    Jihyun%An%NULL%0,                    Jihyun%An%NULL%0]'
$ws.Range("I32").Value = '_PMC_CROSSREF'

$ws.Range("E33").Value = '[Chunyang%Li%NULL%1,                    Fang%Ji%NULL%1,                    Liang%Wang%NULL%1,                    Liping%Wang%NULL%1,                    Jungui%Hao%NULL%1,                    Mingjia%Dai%NULL%1,                    Yan%Liu%NULL%0,                    Xiucheng%Pan%NULL%1,                    Juanjuan%Fu%NULL%1,                    Li%Li%NULL%0,                    Guangde%Yang%NULL%1,                    Jianye%Yang%NULL%1,                    Xuebing%Yan%NULL%1,                    Bing%Gu%NULL%1]'
$ws.Range("I33").Value = '_PMC_CROSSREF'

$ws.Range("E34").Value = '[Diangeng%Li%NULL%1,                    Meiling%Jin%NULL%1,                    Pengtao%Bao%NULL%1,                    Weiguo%Zhao%NULL%1,                    Shixi%Zhang%NULL%1]'
$ws.Range("I34").Value = '_PMC_CROSSREF'

$ws.Range("E35").Value = '[Na%Li%NULL%0,                    Lefei%Han%NULL%2,                    Min%Peng%2658706528@qq.com%2,                    Yuxia%Lv%NULL%2,                    Yin%Ouyang%NULL%2,                    Kui%Liu%NULL%0,                    Linli%Yue%NULL%2,                    Qiannan%Li%NULL%2,                    Guoqiang%Sun%NULL%2,                    Lin%Chen%NULL%2,                    Lin%Yang%l.yang@polyu.edu.hk%0]'
$ws.Range("I35").Value = '_PMC_CROSSREF'

$ws.Range("E36").Value = '[Wei%Li%NULL%1,                    Bo%Zhang%NULL%2,                    Jianhua%Lu%NULL%1,                    Shihua%Liu%NULL%1,                    Zhiqiang%Chang%NULL%1,                    Peng%Cao%NULL%1,                    Xinhua%Liu%NULL%1,                    Peng%Zhang%NULL%2,                    Yan%Ling%NULL%1,                    Kaixiong%Tao%NULL%1,                    Jianying%Chen%Bobytail@sina.com%1]'
$ws.Range("I36").Value = '_PMC_CROSSREF'

$ws.Range("E37").Value = '[Youjiang%Li%NULL%0,                    Yingying%Hu%NULL%2,                    Yingying%Hu%NULL%0,                    Yuanyuan%Yu%NULL%1,                    Xiaodong%Zhang%NULL%1,                    Bin%Li%NULL%1,                    Jianguo%Wu%NULL%1,                    Junyu%Li%NULL%1,                    Yingping%Wu%NULL%1,                    Xiaoping%Xia%NULL%2,                    Xiaoping%Xia%NULL%0,                    Huina%Tang%NULL%1,                    Jian%Xu%xuj@zju.edu.cn%1]'
$ws.Range("I37").Value = '_PMC_CROSSREF'

$ws.Range("E38").Value = '[ Y.%Li%null%1,                     H.% Qian%null%1,                     J.% Hang%null%1,                     X.% Chen%null%1,                     L.% Hong%null%1,                     P.% Liang%null%1,                     J.% Li%null%1,                     S.% Xiao%null%1,                     J.% Wei%null%1,                     L.% Liu%null%2,                     M. % Kang%null%1]'
$ws.Range("I38").Value = '_MedBiorxiv_CROSSREF'

$ws.Range("E39").Value = '[Wei%Liu%NULL%0,                    Jing%Wang%NULL%3,                    Wenbin%Li%NULL%1,                    Zhaoxian%Zhou%NULL%1,                    Siying%Liu%NULL%1,                    Zhihui%Rong%rongzhihui53@163.com%1]'
$ws.Range("I39").Value = '_PMC_Springer_CROSSREF'
